$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order the header columns in row 2 so that the "business key" / id
# columns come first, followed by EndDate and StartDate.
$ws.Range("A2").Value = "OrganizationPersonRole_ID"
$ws.Range("B2").Value = "OrganizationBusinessKey"
$ws.Range("C2").Value = "PersonBusinessKey"
$ws.Range("D2").Value = "RoleBusinessKey"
$ws.Range("E2").Value = "EndDate"
$ws.Range("F2").Value = "StartDate"
